$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H93").Value = 43925
$ws_ALC.Range("J93").Value = 43925
$ws_ALC.Range("L93").Value = 43925
$ws_ALC.Range("N93").Value = -48917
$ws_ALC.Range("H112").Value = 1290.5
$ws_ALC.Range("I112").Value = 450
$ws_ALC.Range("J112").Value = 1334.7368
$ws_ALC.Range("K112").Value = 1350
$ws_ALC.Range("L112").Value = 4004.2104
$ws_ALC.Range("M112").Value = -242
$ws_ALC.Range("N112").Value = -6220.2104
$ws_ALC.Range("H125").Value = 688.13336
$ws_ALC.Range("I125").Value = 515.25
$ws_ALC.Range("K125").Value = 4637.25
$ws_ALC.Range("M125").Value = -2177.25
$ws_ALC.Range("H129").Value = 1314.671
$ws_ALC.Range("I129").Value = 382.85715
$ws_ALC.Range("J129").Value = 1409.2029
$ws_ALC.Range("K129").Value = 1148.57145
$ws_ALC.Range("L129").Value = 4227.6087
$ws_ALC.Range("M129").Value = 3851.42855
$ws_ALC.Range("N129").Value = -14227.6087
$ws_ALC.Range("H132").Value = 26258458
$ws_ALC.Range("I132").Value = 27407160
$ws_ALC.Range("K132").Value = 82221480
$ws_ALC.Range("M132").Value = -82218950
$ws_ALC.Range("H137").Value = 683243.0600000001
$ws_ALC.Range("I137").Value = 1590379
$ws_ALC.Range("J137").Value = 2891.05
$ws_ALC.Range("K137").Value = 4771137
$ws_ALC.Range("L137").Value = 8673.150000000001
$ws_ALC.Range("M137").Value = -4768587
$ws_ALC.Range("N137").Value = -13773.15
$ws_ALC.Range("H138").Value = 2984.6667
$ws_ALC.Range("I138").Value = 2181.5
$ws_ALC.Range("J138").Value = 4055.5557
$ws_ALC.Range("K138").Value = 6544.5
$ws_ALC.Range("L138").Value = 12166.6671
$ws_ALC.Range("M138").Value = -1404.5
$ws_ALC.Range("N138").Value = -22446.6671
$ws_ALC.Range("H141").Value = 31205.285
$ws_ALC.Range("I141").Value = 43532.707
$ws_ALC.Range("K141").Value = 130598.121
$ws_ALC.Range("M141").Value = -125418.121

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H32").Value = 4759.212
$ws_ARM.Range("I32").Value = 4845.566
$ws_ARM.Range("J32").Value = 4407.154
$ws_ARM.Range("K32").Value = 4845.566
$ws_ARM.Range("L32").Value = 4407.154
$ws_ARM.Range("M32").Value = -4558.566
$ws_ARM.Range("N32").Value = -4981.154
$ws_ARM.Range("H45").Value = 1534.25
$ws_ARM.Range("I45").Value = 1478.1875
$ws_ARM.Range("J45").Value = 1646.375
$ws_ARM.Range("K45").Value = 1478.1875
$ws_ARM.Range("L45").Value = 1646.375
$ws_ARM.Range("M45").Value = -1101.1875
$ws_ARM.Range("N45").Value = -2400.375
$ws_ARM.Range("H61").Value = 2965.375
$ws_ARM.Range("I61").Value = 2965.375
$ws_ARM.Range("K61").Value = 2965.375
$ws_ARM.Range("M61").Value = -2753.375
$ws_ARM.Range("H63").Value = 8660350
$ws_ARM.Range("I63").Value = 15391533
$ws_ARM.Range("J63").Value = 5971.4287
$ws_ARM.Range("K63").Value = 15391533
$ws_ARM.Range("L63").Value = 5971.4287
$ws_ARM.Range("M63").Value = -15390847
$ws_ARM.Range("N63").Value = -7343.4287
$ws_ARM.Range("H66").Value = 8660350
$ws_ARM.Range("I66").Value = 15391533
$ws_ARM.Range("J66").Value = 5971.4287
$ws_ARM.Range("K66").Value = 76957665
$ws_ARM.Range("L66").Value = 29857.1435
$ws_ARM.Range("M66").Value = -76954233
$ws_ARM.Range("N66").Value = -36721.14350000001
$ws_ARM.Range("H110").Value = 2752.5
$ws_ARM.Range("I110").Value = 5761
$ws_ARM.Range("J110").Value = 1248.25
$ws_ARM.Range("K110").Value = 5761
$ws_ARM.Range("L110").Value = 1248.25
$ws_ARM.Range("M110").Value = -3716
$ws_ARM.Range("N110").Value = -5338.25
$ws_ARM.Range("H112").Value = 35631.58
$ws_ARM.Range("J112").Value = 35631.58
$ws_ARM.Range("L112").Value = 35631.58
$ws_ARM.Range("N112").Value = -38585.58
$ws_ARM.Range("H122").Value = 3394.0557
$ws_ARM.Range("I122").Value = 3205.12
$ws_ARM.Range("J122").Value = 3823.4546
$ws_ARM.Range("K122").Value = 9615.360000000001
$ws_ARM.Range("L122").Value = 11470.3638
$ws_ARM.Range("M122").Value = -7165.360000000001
$ws_ARM.Range("N122").Value = -16370.3638
$ws_ARM.Range("H132").Value = 2820.1538
$ws_ARM.Range("I132").Value = 2352.4333
$ws_ARM.Range("J132").Value = 4379.222
$ws_ARM.Range("K132").Value = 7057.2999
$ws_ARM.Range("L132").Value = 13137.666
$ws_ARM.Range("M132").Value = -4527.2999
$ws_ARM.Range("N132").Value = -18197.666
$ws_ARM.Range("H136").Value = 2965.375
$ws_ARM.Range("I136").Value = 2965.375
$ws_ARM.Range("K136").Value = 8896.125
$ws_ARM.Range("M136").Value = -6346.125

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H95").Value = 31777.777
$ws_BSM.Range("J95").Value = 31777.777
$ws_BSM.Range("L95").Value = 31777.777
$ws_BSM.Range("N95").Value = -37269.777
$ws_BSM.Range("H107").Value = 1347.2632
$ws_BSM.Range("J107").Value = 1369.6666
$ws_BSM.Range("L107").Value = 1369.6666
$ws_BSM.Range("N107").Value = -5209.6666
$ws_BSM.Range("H134").Value = 3519.65
$ws_BSM.Range("I134").Value = 1277.2667
$ws_BSM.Range("J134").Value = 4865.08
$ws_BSM.Range("K134").Value = 3831.800099999999
$ws_BSM.Range("L134").Value = 14595.24
$ws_BSM.Range("M134").Value = -1296.800099999999
$ws_BSM.Range("N134").Value = -19665.24

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H86").Value = 1904.6897
$ws_CRP.Range("I86").Value = 1755.7368
$ws_CRP.Range("J86").Value = 2187.7
$ws_CRP.Range("K86").Value = 1755.7368
$ws_CRP.Range("L86").Value = 2187.7
$ws_CRP.Range("M86").Value = -632.7367999999999
$ws_CRP.Range("N86").Value = -4433.7
$ws_CRP.Range("H89").Value = 1904.6897
$ws_CRP.Range("I89").Value = 1755.7368
$ws_CRP.Range("J89").Value = 2187.7
$ws_CRP.Range("K89").Value = 8778.683999999999
$ws_CRP.Range("L89").Value = 10938.5
$ws_CRP.Range("M89").Value = -3162.683999999999
$ws_CRP.Range("N89").Value = -22170.5
$ws_CRP.Range("H132").Value = 3026.6365
$ws_CRP.Range("I132").Value = 1641.1875
$ws_CRP.Range("J132").Value = 6721.1665
$ws_CRP.Range("K132").Value = 4923.5625
$ws_CRP.Range("L132").Value = 20163.4995
$ws_CRP.Range("M132").Value = -2393.5625
$ws_CRP.Range("N132").Value = -25223.4995

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H6").Value = 113.875
$ws_CUL.Range("I6").Value = 113.875
$ws_CUL.Range("K6").Value = 341.625
$ws_CUL.Range("M6").Value = -228.625
$ws_CUL.Range("H68").Value = 4708.515
$ws_CUL.Range("J68").Value = 7989.9414
$ws_CUL.Range("L68").Value = 23969.8242
$ws_CUL.Range("N68").Value = -25591.8242
$ws_CUL.Range("H71").Value = 4708.515
$ws_CUL.Range("J71").Value = 7989.9414
$ws_CUL.Range("L71").Value = 71909.47259999999
$ws_CUL.Range("N71").Value = -80021.47259999999
$ws_CUL.Range("H86").Value = 928.34784
$ws_CUL.Range("I86").Value = 503.84616
$ws_CUL.Range("J86").Value = 1480.2
$ws_CUL.Range("K86").Value = 1511.53848
$ws_CUL.Range("L86").Value = 4440.6
$ws_CUL.Range("M86").Value = -325.5384799999999
$ws_CUL.Range("N86").Value = -6812.6
$ws_CUL.Range("H89").Value = 928.34784
$ws_CUL.Range("I89").Value = 503.84616
$ws_CUL.Range("J89").Value = 1480.2
$ws_CUL.Range("K89").Value = 4534.61544
$ws_CUL.Range("L89").Value = 13321.8
$ws_CUL.Range("M89").Value = 1393.38456
$ws_CUL.Range("N89").Value = -25177.8
$ws_CUL.Range("H107").Value = 14485.946
$ws_CUL.Range("I107").Value = 391.55
$ws_CUL.Range("J107").Value = 31067.588
$ws_CUL.Range("K107").Value = 1174.65
$ws_CUL.Range("L107").Value = 93202.764
$ws_CUL.Range("M107").Value = 745.3499999999999
$ws_CUL.Range("N107").Value = -97042.764
$ws_CUL.Range("H113").Value = 1812280.5
$ws_CUL.Range("I113").Value = 624.9815
$ws_CUL.Range("J113").Value = 8334240.5
$ws_CUL.Range("K113").Value = 1874.9445
$ws_CUL.Range("L113").Value = 25002721.5
$ws_CUL.Range("M113").Value = 295.0554999999999
$ws_CUL.Range("N113").Value = -25007061.5
$ws_CUL.Range("H120").Value = 4000
$ws_CUL.Range("I120").Value = 3000
$ws_CUL.Range("J120").Value = 5000
$ws_CUL.Range("K120").Value = 9000
$ws_CUL.Range("L120").Value = 15000
$ws_CUL.Range("M120").Value = -4162
$ws_CUL.Range("N120").Value = -24676
$ws_CUL.Range("H131").Value = 825.02
$ws_CUL.Range("J131").Value = 843.9149
$ws_CUL.Range("L131").Value = 2531.7447
$ws_CUL.Range("N131").Value = -12611.7447

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H113").Value = 1660.4
$ws_GSM.Range("I113").Value = 1827.2858
$ws_GSM.Range("K113").Value = 1827.2858
$ws_GSM.Range("M113").Value = 342.7141999999999
$ws_GSM.Range("H132").Value = 2757.465
$ws_GSM.Range("I132").Value = 2258.742
$ws_GSM.Range("J132").Value = 4045.8333
$ws_GSM.Range("K132").Value = 6776.226000000001
$ws_GSM.Range("L132").Value = 12137.4999
$ws_GSM.Range("M132").Value = -4246.226000000001
$ws_GSM.Range("N132").Value = -17197.4999

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H40").Value = 5229.9287
$ws_LTW.Range("I40").Value = 4018.25
$ws_LTW.Range("K40").Value = 4018.25
$ws_LTW.Range("M40").Value = -3882.25
$ws_LTW.Range("H132").Value = 3797.8518
$ws_LTW.Range("I132").Value = 2724.5454
$ws_LTW.Range("J132").Value = 8520.4
$ws_LTW.Range("K132").Value = 8173.6362
$ws_LTW.Range("L132").Value = 25561.2
$ws_LTW.Range("M132").Value = -5643.6362
$ws_LTW.Range("N132").Value = -30621.2

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H136").Value = 2651.311
$ws_WVR.Range("I136").Value = 1175.68
$ws_WVR.Range("K136").Value = 3527.04
$ws_WVR.Range("M136").Value = -977.04
